# Updated cryptos list with GitHub Actions.
# Refreshes the per-coin Price (D) and Volume(1h) (E) columns scraped from
# coinranking.com, and re-syncs three coin pairs whose scrape order flipped
# between runs (rows 8/9, 43/44, 48/49 swap Name/Link/Price/Volume together).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.127.99"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.165.93"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.41"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.92"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.164.76"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.548"
$ws.Range("E9").Value = "  +2.54%  "
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.56"
$ws.Range("E11").Value = "  -10.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.519"
$ws.Range("E12").Value = "  +2.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000267"
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.33"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.684.10"
$ws.Range("E15").Value = "  -1.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.137.66"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.42"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.163.81"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "510.20"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.38"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.728"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.03"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.66"
$ws.Range("E24").Value = "  -3.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.53"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.99"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.06"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("E29").Value = "  +6.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.06"
$ws.Range("E30").Value = "  +7.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.01"
$ws.Range("E31").Value = "  +3.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.95"
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("E34").Value = "  -1.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.49"
$ws.Range("E35").Value = "  -1.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "497.22"
$ws.Range("E36").Value = "  +3.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.65"
$ws.Range("E37").Value = "  -1.16%  "
$ws.Range("E38").Value = "  -2.45%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.128"
$ws.Range("E40").Value = "  +7.88%  "
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0675"
$ws.Range("E42").Value = "  +5.76%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.295"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("E44").Value = "  -5.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.42"
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.820.41"
$ws.Range("E46").Value = "  -4.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.80"
$ws.Range("E47").Value = "  -3.38%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.37"
$ws.Range("E48").Value = "  +2.95%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.02"
$ws.Range("E51").Value = "  +2.79%  "
